$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of Argent (Silver) price data appended at row 73.
# Source data is stored as text in the sheet (like the rest of the rows),
# so force the number format to Text before assigning values.
$row = 73
$rng = $ws.Range("A$row`:J$row")
$rng.NumberFormat = "@"

$ws.Range("A$row").Value = "2025-05-13"
$ws.Range("B$row").Value = "38"
$ws.Range("C$row").Value = "37.2"
$ws.Range("D$row").Value = "0.98"
$ws.Range("E$row").Value = "0.265"
$ws.Range("F$row").Value = "0.09"
$ws.Range("G$row").Value = "5,318"
$ws.Range("H$row").Value = "7,962"
$ws.Range("I$row").Value = "8,012"
$ws.Range("J$row").Value = "7.2476"
